$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy row 37's formatting down into the two new rows (38, 39) so the new
# cells pick up the same style indices (border/fill/font) as the rest of
# the table instead of being left unstyled.
$ws.Range("A37:E37").Copy($ws.Range("A38:E38"))
$ws.Range("A37:E37").Copy($ws.Range("A39:E39"))

# The existing last row (PublishedAPostTimeStampTest) moves from PASS to
# SKIP now that it is no longer the final executed case.
$ws.Cells.Item(37, 5).Value = "SKIP"

# New test case: CommentsTabTimeStampValidationTest
$ws.Cells.Item(38, 1).Value = "CommentsTabTimeStampValidationTest"
$ws.Cells.Item(38, 2).Value = "TBD"
$ws.Cells.Item(38, 3).Value = "Verify that Comments tab comments displayed with timestamp"
$ws.Cells.Item(38, 4).Value = "Y"
$ws.Cells.Item(38, 5).Value = "SKIP"

# New test case: HCRProfileBadgeTest
$ws.Cells.Item(39, 1).Value = "HCRProfileBadgeTest"
$ws.Cells.Item(39, 2).Value = "TBD"
$ws.Cells.Item(39, 3).Value = "Verify that HCR profile having badge along with their name"
$ws.Cells.Item(39, 4).Value = "Y"
$ws.Cells.Item(39, 5).Value = "PASS"

# Match the author's final selection/scroll position in the sheet view.
$ws.Range("C34").Select()
